{"js": "// The source diff for this revision only reorders XML attributes\n// (e.g. `w:val`/`w:themeColor`/`w:themeShade` on <w:color>, the\n// `xmlns:*` declarations on the <w:document> root, `w:pgSz`/`w:pgMar`\n// on <w:sectPr>, and the many `<w:lsdException>`/`<w:style>` entries in\n// styles.xml) into alphabetical order. Every attribute name/value pair\n// present before the change is still present afterwards \u2014 nothing was\n// added, removed, or re-valued. That kind of pure re-serialization\n// (the commit message calls it \"Moving from 2.0.1 to 2.0.2\", i.e. a\n// tooling/library version bump that changed how the OOXML is emitted)\n// carries no document-model change for Word's object model to apply:\n// paragraphs, runs, formatting, fields, and section properties are all\n// identical before and after.\n//\n// Word's JavaScript API only ever mutates the document through its\n// object model (Body/Paragraph/Range/Font/PageSetup/...), and any part\n// it rewrites is re-serialized with the host's own fixed attribute\n// order \u2014 there is no supported call that reorders the literal XML\n// attributes of an OOXML part without also touching its content. Since\n// there is no content change to make, the correct replay of this diff\n// is therefore a no-op: touch nothing, so the package stays exactly as\n// it was.\n", "ps1": "# The source diff for this revision only reorders XML attributes\n# (e.g. `w:val`/`w:themeColor`/`w:themeShade` on <w:color>, the\n# `xmlns:*` declarations on the <w:document> root, `w:pgSz`/`w:pgMar`\n# on <w:sectPr>, and the many `<w:lsdException>`/`<w:style>` entries in\n# styles.xml) into alphabetical order. Every attribute name/value pair\n# present before the change is still present afterwards - nothing was\n# added, removed, or re-valued. That kind of pure re-serialization\n# (the commit message calls it \"Moving from 2.0.1 to 2.0.2\", i.e. a\n# tooling/library version bump that changed how the OOXML is emitted)\n# carries no document-model change for Word to apply: paragraphs,\n# runs, formatting, fields, and section properties are all identical\n# before and after.\n#\n# Word's COM object model only mutates the document through its\n# object model (Paragraphs/Range/Find/Font/PageSetup/...), and any\n# part it rewrites is re-serialized with the host's own fixed\n# attribute order - e.g. $d.WordOpenXML is read-only here, and\n# Range.InsertXML reparses/replaces the addressed range's content\n# rather than preserving a caller-chosen literal attribute order.\n# There is no supported call that reorders the literal XML attributes\n# of an OOXML part without also touching its content. Since there is\n# no content change to make, the correct replay of this diff is\n# therefore a no-op: touch nothing, so the document stays exactly as\n# it was.\n"}
